# "Add files via upload" -- appends 10 new crypto-pair entries to the
# existing list, inserts one more ("cryptoAixbtBase.xlsx") right after
# "cryptoAerodromeBase.xlsx", and paints the whole (now 18-row) list with
# a light-green/olive highlight fill, leaving the final two helper rows
# (19-20) on the plain/no-fill style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final column-A contents, row 1 .. row 20, in order.
$values = @(
    "cryptoAAVEpolygon.xlsx",
    "cryptoAerodromeBase.xlsx",
    "cryptoAixbtBase.xlsx",
    "cryptoChainlinkPolygon.xlsx",
    "cryptoDogeBnb.xlsx",
    "cryptoMorphoBase.xlsx",
    "cryptoMystPolygon.xlsx",
    "cryptoPaxgoldPolygon.xlsx",
    "cryptoSolanaPolygon.xlsx",
    "cryptoWrappedBTCPolygon.xlsx",
    "cryptoUniswapPolygon.xlsx",
    "cryptoAtomcosmosBnb.xlsx",
    "cryptoNexoPolygon.xlsx",
    "cryptoVirtualBase.xlsx",
    "cryptoPancakeswapBase.xlsx",
    "cryptoWrappedEthPolygon.xlsx",
    "cryptoWrappedbnbBnb.xlsx",
    "cryptoTelcoinPolygon.xlsx",
    "cryptoXrpBnb.xlsx",
    "cryptoAsterBnb.xlsx"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Rows 1-18 get the new highlight fill (solid, RGB D4EA6B); rows 19-20
# are left on the workbook's plain style.
$ws.Range("A1:A18").Interior.Color = 7072468

# Selection finishes on the last entered cell, matching the saved view.
$ws.Range("A20").Select()
